# Update the "Habitat" domain description (row 24, column C) with the new
# text describing the domain, replacing the previous description.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eixos")

$ws.Range("C24").Value = "Este domínio engloba a fileira do Habitat, focando-se na inovação e sustentabilidade em materiais, produtos e soluções aplicadas à construção, mobiliário, arquitetura e design. Inclui também o desenvolvimento de novos materiais e processos produtivos eficientes, com forte ligação a setores como cortiça, cerâmica, madeira, vidro, metais, domótica e têxteis-lar. Valoriza a eco-inovação, a eficiência energética, o design sustentável e a internacionalização de soluções habitacionais e urbanas."

# Move the active selection from D29 to D30, matching where the author's
# cursor ended up after the edit.
$ws.Range("D30").Select()
